# Rename the embedded drawing parts' docPr/cNvPr "name" attributes:
#   word/footer1.xml : id="3" / id="0"  image2.png -> image1.png  (Pearson logo)
#   word/footer2.xml : id="2" / id="0"  image2.png -> image1.png  (Pearson logo)
#   word/header1.xml : id="1" / id="0"  image1.jpg -> image2.jpg  (BTEC logo)
#
# InlineShape has no settable "Name" property in the Word object model
# (only the floating Shape object does), so the rename is applied by
# round-tripping the document's Open XML package text and swapping the
# handful of distinguishing `id="N" name="..."` tokens.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$xml = $xml.Replace('id="3" name="image2.png"', 'id="3" name="image1.png"')
$xml = $xml.Replace('id="2" name="image2.png"', 'id="2" name="image1.png"')
$xml = $xml.Replace('id="0" name="image2.png"', 'id="0" name="image1.png"')

$xml = $xml.Replace('id="1" name="image1.jpg"', 'id="1" name="image2.jpg"')
$xml = $xml.Replace('id="0" name="image1.jpg"', 'id="0" name="image2.jpg"')

$d.WordOpenXML = $xml
